$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Investments")

# Adjust max weight on Real Estate (IYR, row 7) from 0.1 to 0.05
$ws.Range("C7").Value = 0.05

# Update the active cell selection to match the saved view state
$ws.Range("G12").Select()
